# Auto-update Data Telemetría - Ejecución Diaria
# Appends the new day's (45997) connectivity stats to the "historico" sheet
# and refreshes the "ultimo_snapshot" sheet with that same latest data.

$wb = $excel.ActiveWorkbook

$historico = $wb.Worksheets.Item("historico")
$snapshot  = $wb.Worksheets.Item("ultimo_snapshot")

# New rows for 2025-12-06 (Excel serial date 45997), in column order:
# fecha, resumen, total_vin, cnt_Conectado 0-2, cnt_Intermitente 3-14,
# cnt_Limitado 15-30+, cnt_Desconectado 31+, cnt_Nunca, pct_Conectado 0-2,
# pct_Intermitente 3-14, pct_Limitado 15-30+, pct_Desconectado 31+, pct_Nunca
$newRows = @(
    @(45997, "Telemetría", 5905, 3538, 506, 190, 663, 1008, 59.92, 8.57, 3.22, 11.23, 17.07),
    @(45997, "GPS (según REGLA)", 5305, 4682, 346, 92, 178, 7, 88.26000000000001, 6.52, 1.73, 3.36, 0.13),
    @(45997, "GPS (todas con gps_timestamp)", 11203, 9539, 798, 289, 577, 0, 85.15000000000001, 7.12, 2.58, 5.15, 0)
)

# Append the three new rows to the bottom of "historico" (rows 59-61).
$startRow = $historico.UsedRange.Rows.Count + 1
for ($i = 0; $i -lt $newRows.Length; $i++) {
    $row = $startRow + $i
    $values = $newRows[$i]

    $dateCell = $historico.Cells.Item($row, 1)
    $dateCell.Value = $values[0]
    $dateCell.NumberFormat = $historico.Cells.Item($row - 1, 1).NumberFormat

    for ($col = 2; $col -le $values.Length; $col++) {
        $historico.Cells.Item($row, $col).Value = $values[$col - 1]
    }
}

# Overwrite the "ultimo_snapshot" sheet (rows 2-4) with the same latest data.
for ($i = 0; $i -lt $newRows.Length; $i++) {
    $row = 2 + $i
    $values = $newRows[$i]

    $dateCell = $snapshot.Cells.Item($row, 1)
    $dateCell.Value = $values[0]
    $dateCell.NumberFormat = $snapshot.Cells.Item($row, 1).NumberFormat

    for ($col = 2; $col -le $values.Length; $col++) {
        $snapshot.Cells.Item($row, $col).Value = $values[$col - 1]
    }
}
